$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Genre table (A7:B10): rename the "genre" column to "name"
# -------------------------------------------------------------------------
$ws.Range("B8").Value = "name     charvar(20)"

# -------------------------------------------------------------------------
# User table (H1:O5) gains a new "duration" lead-in column and an "age"
# column (replacing "age_range"), and a brand new "Profile" side-table is
# introduced to its right (Q1:W5). Net effect: everything from H2:N2
# shifts one column to the right (G2:M2), country moves out to the new
# Profile table, and "age_range integer" is replaced with "age integer".
# -------------------------------------------------------------------------

# New "duration" column, then shift the rest of the User columns right by one
$ws.Range("H2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Value = "duration     time"
$ws.Range("H2").Clear() | Out-Null
$ws.Range("I2").Value = "user_id     integer"
$ws.Range("J2").Value = "Name      charvar(35)"
$ws.Range("K2").Value = "email      text"
$ws.Range("L2").Value = "password text"
$ws.Range("M2").Value = "gender    char(1)"
$ws.Range("N2").Value = "age integer"
$ws.Range("O2").Value = "Picture    image"

# Formatting + borders for the new/shifted cells in rows 3-5 under "duration"
$ws.Range("H3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Clear() | Out-Null
$ws.Range("H4").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Clear() | Out-Null
$ws.Range("F5").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null

# New "Profile" side-table: header + columns + blank bordered rows
$ws.Range("H1").Copy() | Out-Null
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q1").Value = "Profile"

$ws.Range("I2").Copy() | Out-Null
$ws.Range("R2").PasteSpecial(-4122) | Out-Null
$ws.Range("R2").Value = "profile_id    integer"
$ws.Range("S2").PasteSpecial(-4122) | Out-Null
$ws.Range("S2").Value = "province    text"
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$ws.Range("T2").Value = "city              text"
$ws.Range("U2").PasteSpecial(-4122) | Out-Null
$ws.Range("U2").Value = "occupation  charvar(30)"
$ws.Range("V2").PasteSpecial(-4122) | Out-Null
$ws.Range("V2").Value = "device-used charvar(25)"
$ws.Range("W2").PasteSpecial(-4122) | Out-Null
$ws.Range("W2").Value = "country   text "

$ws.Range("G3").Copy() | Out-Null
$ws.Range("R3:W3").PasteSpecial(-4122) | Out-Null
$ws.Range("R4:W4").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------------
# Watched / Wish tables (G7:J10, L7:O10): "Wish" header shifts from L7 to
# M7, and both tables lose their trailing "rating" column (J8/J9/J10 and
# O8/O9/O10).
# -------------------------------------------------------------------------
$ws.Range("L7").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Value = "Wish"
$ws.Range("L7").Clear() | Out-Null

$ws.Range("J8").Clear() | Out-Null
$ws.Range("J9").Clear() | Out-Null
$ws.Range("J10").Clear() | Out-Null
$ws.Range("O8").Clear() | Out-Null
$ws.Range("O9").Clear() | Out-Null
$ws.Range("O10").Clear() | Out-Null

# -------------------------------------------------------------------------
# "Role" table (A17:C20) becomes "Plays in", whose first column is now
# "movie_id" instead of "role_id".
# -------------------------------------------------------------------------
$ws.Range("A17").Value = "Plays in"
$ws.Range("A18").Value = "movie_id   integer"

# -------------------------------------------------------------------------
# Actorplays table (A22:K25) gains a trailing "Time_stamp" column.
# -------------------------------------------------------------------------
$ws.Range("K23").Copy() | Out-Null
$ws.Range("L23").PasteSpecial(-4122) | Out-Null
$ws.Range("L23").Value = "Time_stamp time"
$ws.Range("K24").Copy() | Out-Null
$ws.Range("L24").PasteSpecial(-4122) | Out-Null
$ws.Range("K25").Copy() | Out-Null
$ws.Range("L25").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------------
# Profile table (A32:E35) is removed from the bottom of the sheet (it now
# lives next to "User" at the top). "UserLikes Actor" (previously at
# G32:H35) takes its place at A32:B35.
# -------------------------------------------------------------------------
$ws.Range("A32").Value = "UserLikes Actor"
$ws.Range("G32").Clear() | Out-Null

$ws.Range("A33").Value = "user_id     integer"
$ws.Range("B33").Value = "actor_id     integer"
$ws.Range("C33").Clear() | Out-Null
$ws.Range("D33").Clear() | Out-Null
$ws.Range("E33").Clear() | Out-Null
$ws.Range("G33").Clear() | Out-Null
$ws.Range("H33").Clear() | Out-Null

$ws.Range("C34").Clear() | Out-Null
$ws.Range("D34").Clear() | Out-Null
$ws.Range("E34").Clear() | Out-Null
$ws.Range("G34").Clear() | Out-Null
$ws.Range("H34").Clear() | Out-Null
$ws.Range("C35").Clear() | Out-Null
$ws.Range("D35").Clear() | Out-Null
$ws.Range("E35").Clear() | Out-Null
$ws.Range("G35").Clear() | Out-Null
$ws.Range("H35").Clear() | Out-Null

# -------------------------------------------------------------------------
# Cosmetic sheet-level adjustments to match the final view/column sizing.
# -------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 57.6
$ws.Columns.Item(11).ColumnWidth = 10.17
$ws.Columns.Item(12).ColumnWidth = 10.5
$ws.Columns.Item(17).ColumnWidth = 10.17

$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.Zoom = 100
$ws.Range("J18").Select() | Out-Null
